# Update crypto price/volume data per the source feed refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text value (matches the original
    # inlineStr cell type) even when the new string looks numeric,
    # then drop the temporary "@" number format so the cell keeps
    # its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '27.707.33'
$ws.Range("E2").Value = '  +2.96%  '

Set-TextValue $ws.Range("D3") '1.862.15'
$ws.Range("E3").Value = '  +2.66%  '

Set-TextValue $ws.Range("D4") '1.037'
$ws.Range("E4").Value = '  +3.11%  '

Set-TextValue $ws.Range("D5") '324.27'
$ws.Range("E5").Value = '  +3.95%  '

Set-TextValue $ws.Range("D6") '1.032'
$ws.Range("E6").Value = '  +2.72%  '

Set-TextValue $ws.Range("D7") '0.4401'
$ws.Range("E7").Value = '  +2.51%  '

Set-TextValue $ws.Range("D8") '0.3790'
$ws.Range("E8").Value = '  +2.44%  '

$ws.Range("E9").Value = '  +3.04%  '

Set-TextValue $ws.Range("D10") '0.8831'
$ws.Range("E10").Value = '  +1.77%  '

$ws.Range("E11").Value = '  +1.97%  '

Set-TextValue $ws.Range("D12") '1.882.60'
$ws.Range("E12").Value = '  -9.25%  '

$ws.Range("E13").Value = '  +2.81%  '

Set-TextValue $ws.Range("D14") '6.741'
$ws.Range("E14").Value = '  +1.61%  '

Set-TextValue $ws.Range("D15") '0.07214'
$ws.Range("E15").Value = '  +4.06%  '

Set-TextValue $ws.Range("D16") '83.73'

$ws.Range("E17").Value = '  +3.28%  '

$ws.Range("E18").Value = '  +2.92%  '

$ws.Range("E19").Value = '  +2.77%  '

$ws.Range("E20").Value = '  +2.07%  '

Set-TextValue $ws.Range("D21") '27.729.51'
$ws.Range("E21").Value = '  +2.95%  '

Set-TextValue $ws.Range("D22") '5.303'
$ws.Range("E22").Value = '  +1.95%  '

Set-TextValue $ws.Range("D23") '11.39'
$ws.Range("E23").Value = '  +3.77%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D24") '1.955'
$ws.Range("E24").Value = '  +3.65%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D25") '158.06'
$ws.Range("E25").Value = '  +2.27%  '

Set-TextValue $ws.Range("D26") '18.83'
$ws.Range("E26").Value = '  +2.53%  '

Set-TextValue $ws.Range("D27") '1.999'
$ws.Range("E27").Value = '  +3.60%  '

Set-TextValue $ws.Range("D28") '5.299'
$ws.Range("E28").Value = '  +0.94%  '

Set-TextValue $ws.Range("D29") '117.34'
$ws.Range("E29").Value = '  +2.42%  '

Set-TextValue $ws.Range("D30") '0.09079'
$ws.Range("E30").Value = '  +1.38%  '

Set-TextValue $ws.Range("D31") '1.213'
$ws.Range("E31").Value = '  +3.80%  '

Set-TextValue $ws.Range("D32") '0.7714'
$ws.Range("E32").Value = '  +3.58%  '

Set-TextValue $ws.Range("D33") '3.030'
$ws.Range("E33").Value = '  +8.14%  '

Set-TextValue $ws.Range("D34") '4.576'
$ws.Range("E34").Value = '  +3.26%  '

Set-TextValue $ws.Range("D35") '1.034'
$ws.Range("E35").Value = '  +2.91%  '

Set-TextValue $ws.Range("D36") '1.165'
$ws.Range("E36").Value = '  +4.53%  '

$ws.Range("E37").Value = '  +3.36%  '

Set-TextValue $ws.Range("D38") '0.05345'
$ws.Range("E38").Value = '  +2.00%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D39") '0.5188'
$ws.Range("E39").Value = '  +1.67%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D40") '2.836'
$ws.Range("E40").Value = '  +3.03%  '

Set-TextValue $ws.Range("D41") '0.1690'
$ws.Range("E41").Value = '  +2.09%  '

Set-TextValue $ws.Range("D42") '6.844'
$ws.Range("E42").Value = '  +5.11%  '

Set-TextValue $ws.Range("D43") '8.704'
$ws.Range("E43").Value = '  +4.68%  '

Set-TextValue $ws.Range("D44") '109.46'
$ws.Range("E44").Value = '  +1.99%  '

Set-TextValue $ws.Range("D45") '10.69'
$ws.Range("E45").Value = '  +2.32%  '

Set-TextValue $ws.Range("D46") '1.720'
$ws.Range("E46").Value = '  +4.26%  '

Set-TextValue $ws.Range("D47") '0.4695'
$ws.Range("E47").Value = '  +2.39%  '

Set-TextValue $ws.Range("D48") '0.06424'
$ws.Range("E48").Value = '  +2.39%  '

Set-TextValue $ws.Range("D49") '1.868'
$ws.Range("E49").Value = '  +3.18%  '

Set-TextValue $ws.Range("D50") '39.82'
$ws.Range("E50").Value = '  +4.30%  '

Set-TextValue $ws.Range("D51") '64.50'
$ws.Range("E51").Value = '  +1.97%  '
